$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 311, shifting rows 311:363 down to 312:364.
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row 311 with its data.
$ws.Cells.Item(311, 1).Value = 8
$ws.Cells.Item(311, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(311, 3).Value = "Coquimbo"
$ws.Cells.Item(311, 4).Value = 44951
$ws.Cells.Item(311, 5).Value = 4
$ws.Cells.Item(311, 6).Value = 100112012
$ws.Cells.Item(311, 7).Value = "Espinaca"
$ws.Cells.Item(311, 8).Value = "Sin especificar"
$ws.Cells.Item(311, 9).Value = "Primera"
$ws.Cells.Item(311, 10).Value = 2000
$ws.Cells.Item(311, 11).Value = 500
$ws.Cells.Item(311, 12).Value = 600
$ws.Cells.Item(311, 13).Value = 550
$ws.Cells.Item(311, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(311, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(311, 16).Value = 1100
$ws.Cells.Item(311, 17).Value = 0.5
$ws.Cells.Item(311, 18).Value = "Hortaliza"
